$d = $word.ActiveDocument
$t = $d.Tables(1)
$vbreak = [char]11

$t.Cell(1,1).Range.Text = "61 x 85" + $vbreak + "  8    5" + $vbreak + "  ----" + $vbreak + "6|    |" + $vbreak + "1|    |"
$t.Cell(1,2).Range.Text = "65 x 96" + $vbreak + "  9    6" + $vbreak + "  ----" + $vbreak + "6|    |" + $vbreak + "5|    |"
$t.Cell(1,3).Range.Text = "43 x 98" + $vbreak + "  9    8" + $vbreak + "  ----" + $vbreak + "4|    |" + $vbreak + "3|    |"

$t.Cell(2,1).Range.Text = "94 x 59" + $vbreak + "  5    9" + $vbreak + "  ----" + $vbreak + "9|    |" + $vbreak + "4|    |"
$t.Cell(2,2).Range.Text = "87 x 36" + $vbreak + "  3    6" + $vbreak + "  ----" + $vbreak + "8|    |" + $vbreak + "7|    |"
$t.Cell(2,3).Range.Text = "62 x 44" + $vbreak + "  4    4" + $vbreak + "  ----" + $vbreak + "6|    |" + $vbreak + "2|    |"

$t.Cell(3,1).Range.Text = "88 x 80" + $vbreak + "  8    0" + $vbreak + "  ----" + $vbreak + "8|    |" + $vbreak + "8|    |"
$t.Cell(3,2).Range.Text = "44 x 15" + $vbreak + "  1    5" + $vbreak + "  ----" + $vbreak + "4|    |" + $vbreak + "4|    |"
$t.Cell(3,3).Range.Text = "15 x 32" + $vbreak + "  3    2" + $vbreak + "  ----" + $vbreak + "1|    |" + $vbreak + "5|    |"

$t.Cell(4,1).Range.Text = "69 x 71" + $vbreak + "  7    1" + $vbreak + "  ----" + $vbreak + "6|    |" + $vbreak + "9|    |"
$t.Cell(4,2).Range.Text = "47 x 81" + $vbreak + "  8    1" + $vbreak + "  ----" + $vbreak + "4|    |" + $vbreak + "7|    |"
$t.Cell(4,3).Range.Text = "20 x 39" + $vbreak + "  3    9" + $vbreak + "  ----" + $vbreak + "2|    |" + $vbreak + "0|    |"

$t.Cell(5,1).Range.Text = "94 x 28" + $vbreak + "  2    8" + $vbreak + "  ----" + $vbreak + "9|    |" + $vbreak + "4|    |"
$t.Cell(5,2).Range.Text = "27 x 36" + $vbreak + "  3    6" + $vbreak + "  ----" + $vbreak + "2|    |" + $vbreak + "7|    |"
$t.Cell(5,3).Range.Text = "42 x 54" + $vbreak + "  5    4" + $vbreak + "  ----" + $vbreak + "4|    |" + $vbreak + "2|    |"

Write-Host "Updated 15 lattice multiplication cells"
